# This script reproduces, via Excel COM automation, the edits described by the
# target diff:
#   - On each of the six data/result worksheets (Summary information,
#     Table 1 Submission, Table 2 Authorizations, Table 3 Actions,
#     Table 4 Holdings, Table 5 Auth. entities) a blank cell A1 is introduced
#     (extending each sheet's dimension to start at column A) and a cell
#     comment with the text "All fields found." is attached to A1.
#   - The "Syntax check results" sheet content is restructured/expanded with
#     a new "Structure check" / "Content check" layout and "Link" hyperlink
#     cells pointing at the relevant sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper data: the six sheets that receive a new A1 cell + "All fields found."
# comment, and (besides "Summary information", whose row 1 already exists)
# how many extra blank leading rows need to be materialized.
# ---------------------------------------------------------------------------

$commentText = "All fields found."

# Summary information: row 1 already has content (B1), just need A1.
$ws = $wb.Worksheets.Item("Summary information")
$ws.Cells.Item(1, 1).Style = "Normal"
$ws.Range("A1").AddComment($commentText)

# Table 1 Submission: needs new rows 1 (A1) and 2 (blank).
$ws = $wb.Worksheets.Item("Table 1 Submission")
$ws.Cells.Item(1, 1).Style = "Normal"
$ws.Rows.Item(2).OutlineLevel = 0
$ws.Range("A1").AddComment($commentText)

# Table 2 Authorizations: needs new rows 1 (A1) and 2 (blank).
$ws = $wb.Worksheets.Item("Table 2 Authorizations")
$ws.Cells.Item(1, 1).Style = "Normal"
$ws.Rows.Item(2).OutlineLevel = 0
$ws.Range("A1").AddComment($commentText)

# Table 3 Actions: needs new rows 1 (A1) and 2 (blank).
$ws = $wb.Worksheets.Item("Table 3 Actions")
$ws.Cells.Item(1, 1).Style = "Normal"
$ws.Rows.Item(2).OutlineLevel = 0
$ws.Range("A1").AddComment($commentText)

# Table 4 Holdings: needs new rows 1 (A1) and 2 (blank).
$ws = $wb.Worksheets.Item("Table 4 Holdings")
$ws.Cells.Item(1, 1).Style = "Normal"
$ws.Rows.Item(2).OutlineLevel = 0
$ws.Range("A1").AddComment($commentText)

# Table 5 Auth. entities: needs new rows 1 (A1) and 2-7 (blank).
$ws = $wb.Worksheets.Item("Table 5 Auth. entities")
$ws.Cells.Item(1, 1).Style = "Normal"
$ws.Rows.Item(2).OutlineLevel = 0
$ws.Rows.Item(3).OutlineLevel = 0
$ws.Rows.Item(4).OutlineLevel = 0
$ws.Rows.Item(5).OutlineLevel = 0
$ws.Rows.Item(6).OutlineLevel = 0
$ws.Rows.Item(7).OutlineLevel = 0
$ws.Range("A1").AddComment($commentText)

# ---------------------------------------------------------------------------
# "Syntax check results" sheet: rebuild the body (everything below A1) with
# the new Structure check / Content check layout.
# ---------------------------------------------------------------------------

$ws = $wb.Worksheets.Item("Syntax check results")

# Wipe out everything but A1 ("Guyana 2022 correct.xlsx"), which is unchanged.
$ws.Range("A2:E30").ClearContents()

$ws.Range("B2").Value = "Structure check"

$ws.Range("C3").Value = "Correct number of worksheets in workbook"
$ws.Range("C4").Value = "All worksheets found in workbook."

$ws.Range("C5").Value = "Summary information: Table 1: Submission"
$ws.Range("D6").Value = "Link"
$ws.Range("E6").Value = "All fields found."

$ws.Range("C7").Value = "Summary information: Table 2: Authorizations"
$ws.Range("D8").Value = "Link"
$ws.Range("E8").Value = "All fields found."

$ws.Range("C9").Value = "Summary information: Table 3: Actions"
$ws.Range("D10").Value = "Link"
$ws.Range("E10").Value = "All fields found."

$ws.Range("C11").Value = "Summary information: Table 4: Holdings"
$ws.Range("D12").Value = "Link"
$ws.Range("E12").Value = "All fields found."

$ws.Range("C13").Value = "Table 1 Submission"
$ws.Range("D14").Value = "Link"
$ws.Range("E14").Value = "All fields found."

$ws.Range("C15").Value = "Table 2 Authorizations"
$ws.Range("D16").Value = "Link"
$ws.Range("E16").Value = "All fields found."

$ws.Range("C17").Value = "Table 3 Actions"
$ws.Range("D18").Value = "Link"
$ws.Range("E18").Value = "All fields found."

$ws.Range("C19").Value = "Table 4 Holdings"
$ws.Range("D20").Value = "Link"
$ws.Range("E20").Value = "All fields found."

$ws.Range("C21").Value = "Table 5 Auth. entities"
$ws.Range("D22").Value = "Link"
$ws.Range("E22").Value = "All fields found."

$ws.Range("B23").Value = "Content check"
$ws.Range("C24").Value = "Table 1 Submission"
$ws.Range("C25").Value = "Table 2 Authorizations"
$ws.Range("C26").Value = "Table 3 Actions"
$ws.Range("C27").Value = "Table 4 Holdings"
$ws.Range("C28").Value = "Table 5 Auth. entities"

$ws.Range("A30").Value = "Syntax check completed successfully."

# Hyperlinks on the "Link" cells, pointing at the relevant sheet/section.
$ws.Hyperlinks.Add($ws.Range("D6"),  "Summary information!B3")
$ws.Hyperlinks.Add($ws.Range("D8"),  "Summary information!B17")
$ws.Hyperlinks.Add($ws.Range("D10"), "Summary information!B40")
$ws.Hyperlinks.Add($ws.Range("D12"), "Summary information!B72")
$ws.Hyperlinks.Add($ws.Range("D14"), "Table 1 Submission!B3")
$ws.Hyperlinks.Add($ws.Range("D16"), "Table 2 Authorizations!B3")
$ws.Hyperlinks.Add($ws.Range("D18"), "Table 3 Actions!B3")
$ws.Hyperlinks.Add($ws.Range("D20"), "Table 4 Holdings!B3")
$ws.Hyperlinks.Add($ws.Range("D22"), "Table 5 Auth. entities!C8")

# Restore plain (non-hyperlink) styling on the "Link" cells so they match the
# rest of the sheet, which carries no explicit cell styling. Each cell is
# reset individually (rather than via one D6:D22 range) so that the rows in
# between (7, 9, 11, ...), which have no cell in column D at all, stay empty.
$ws.Range("D6").Value  = "Link"
$ws.Range("D6").Style  = "Normal"
$ws.Range("D8").Value  = "Link"
$ws.Range("D8").Style  = "Normal"
$ws.Range("D10").Value = "Link"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Value = "Link"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Value = "Link"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Value = "Link"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Value = "Link"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Value = "Link"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Value = "Link"
$ws.Range("D22").Style = "Normal"
